$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 4).Value = '44.987.61'
$ws.Cells.Item(2, 5).Value = '  +0.38%  '

$ws.Cells.Item(3, 4).Value = '2.265.37'
$ws.Cells.Item(3, 5).Value = '  +0.76%  '

$ws.Cells.Item(4, 5).Value = '  -0.63%  '

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '301.99'
$ws.Cells.Item(5, 5).Value = '  -1.39%  '

$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '94.68'
$ws.Cells.Item(6, 5).Value = '  -1.59%  '

$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '0.564'
$ws.Cells.Item(7, 5).Value = '  -1.47%  '

$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = '0.999'
$ws.Cells.Item(8, 5).Value = '  -0.60%  '

$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '0.509'
$ws.Cells.Item(9, 5).Value = '  -2.36%  '

$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '34.27'
$ws.Cells.Item(10, 5).Value = '  -3.30%  '

$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '0.0789'
$ws.Cells.Item(11, 5).Value = '  -2.05%  '

$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '7.19'
$ws.Cells.Item(12, 5).Value = '  -0.56%  '

$ws.Cells.Item(13, 5).Value = '  -1.11%  '

$ws.Cells.Item(14, 4).Value = '2.608.62'
$ws.Cells.Item(14, 5).Value = '  +0.65%  '

$ws.Cells.Item(15, 4).Value = '2.263.09'
$ws.Cells.Item(15, 5).Value = '  +0.83%  '

$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '13.63'
$ws.Cells.Item(16, 5).Value = '  +0.15%  '

$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = '0.799'
$ws.Cells.Item(17, 5).Value = '  -5.11%  '

$ws.Cells.Item(18, 4).Value = '44.850.05'
$ws.Cells.Item(18, 5).Value = '  +0.73%  '

$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '12.94'
$ws.Cells.Item(19, 5).Value = '  +8.20%  '

$ws.Cells.Item(20, 4).Value = '0.0₃0922'
$ws.Cells.Item(20, 5).Value = '  -3.22%  '

$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '6.09'
$ws.Cells.Item(21, 5).Value = '  -3.58%  '

$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '65.52'
$ws.Cells.Item(22, 5).Value = '  +0.02%  '

$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '238.82'
$ws.Cells.Item(23, 5).Value = '  -0.69%  '

$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '2.90'
$ws.Cells.Item(24, 5).Value = '  -2.13%  '

$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '0.998'
$ws.Cells.Item(25, 5).Value = '  -0.43%  '

$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '1.89'
$ws.Cells.Item(26, 5).Value = '  -4.74%  '

$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '41.42'
$ws.Cells.Item(27, 5).Value = '  +10.24%  '

$ws.Cells.Item(28, 5).Value = '  -1.37%  '

$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '9.57'
$ws.Cells.Item(29, 5).Value = '  -3.15%  '

$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '19.56'
$ws.Cells.Item(30, 5).Value = '  -1.95%  '

$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '153.20'
$ws.Cells.Item(31, 5).Value = '  +1.55%  '

$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '5.57'
$ws.Cells.Item(32, 5).Value = '  -7.18%  '

$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '0.0787'
$ws.Cells.Item(33, 5).Value = '  -1.46%  '

$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '2.57'
$ws.Cells.Item(34, 5).Value = '  -2.31%  '

$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '2.97'
$ws.Cells.Item(35, 5).Value = '  -3.83%  '

$ws.Cells.Item(36, 5).Value = '  -1.78%  '

$ws.Cells.Item(37, 5).Value = '  -3.99%  '

$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '1.77'
$ws.Cells.Item(38, 5).Value = '  -4.80%  '

$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '0.0311'
$ws.Cells.Item(39, 5).Value = '  +2.21%  '

$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '3.88'
$ws.Cells.Item(40, 5).Value = '  +1.78%  '

$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '3.25'
$ws.Cells.Item(41, 5).Value = '  -4.86%  '

$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '13.67'
$ws.Cells.Item(42, 5).Value = '  -8.62%  '

$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '0.999'
$ws.Cells.Item(43, 5).Value = '  -0.88%  '

$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '1.94'
$ws.Cells.Item(44, 5).Value = '  +12.41%  '

$ws.Cells.Item(45, 4).Value = '1.752.84'
$ws.Cells.Item(45, 5).Value = '  -4.61%  '

$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '0.197'
$ws.Cells.Item(46, 5).Value = '  +2.96%  '

$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '70.92'
$ws.Cells.Item(47, 5).Value = '  +2.21%  '

$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '75.39'
$ws.Cells.Item(48, 5).Value = '  -5.61%  '

$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '96.42'
$ws.Cells.Item(49, 5).Value = '  -2.71%  '

$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '53.92'
$ws.Cells.Item(50, 5).Value = '  -1.75%  '

$ws.Cells.Item(51, 2).Value = 'THORChain'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '4.66'
$ws.Cells.Item(51, 5).Value = '  -5.15%  '
